$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the Price (D) column cells so numeric-looking
# strings (e.g. "525.00", "0.0000137") are not auto-converted to numbers,
# matching the original inline-string cell content.
$priceCells = "D2","D3","D5","D6","D13","D14","D16","D17","D18","D21","D22","D23","D26","D27","D28","D32","D33","D36","D38","D41","D43","D44","D47","D49","D50"
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "59.006.68"
$ws.Range("E2").Value = "  -2.34%  "
$ws.Range("D3").Value = "2.663.58"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("D5").Value = "525.00"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "144.39"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("E9").Value = "  +8.06%  "
$ws.Range("E10").Value = "  -2.32%  "
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("D13").Value = "3.133.61"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").Value = "59.004.76"
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.675.21"
$ws.Range("E16").Value = "  -3.23%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  -1.74%  "
$ws.Range("D18").Value = "339.20"
$ws.Range("E18").Value = "  -3.31%  "
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("D21").Value = "6.41"
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "64.26"
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value = "0.0₃0802"
$ws.Range("E27").Value = "  -1.74%  "
$ws.Range("D28").Value = "7.11"
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("E29").Value = "  -2.48%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D32").Value = "18.87"
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("D33").Value = "150.57"
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("E34").Value = "  -3.86%  "
$ws.Range("E35").Value = "  -3.98%  "
$ws.Range("D36").Value = "0.895"
$ws.Range("E36").Value = "  -5.80%  "
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").Value = "36.90"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  -5.99%  "
$ws.Range("E40").Value = "  -2.97%  "
$ws.Range("D41").Value = "0.617"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").Value = "276.39"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").Value = "19.90"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("D47").Value = "2.051.57"
$ws.Range("E47").Value = "  -3.63%  "
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0229"
$ws.Range("E49").Value = "  -2.42%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "4.70"
$ws.Range("E50").Value = "  -3.73%  "
$ws.Range("E51").Value = "  -1.62%  "

# Restore default (Normal) style on the price cells so no stray number
# format / quote-prefix style lingers on the cell once the text value is set.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
